$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-e-asia-joint-research-program-grants", "AEDTApplications"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/partnership-projects", "ASAP"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/development-grants", "ASAPRegister"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/ideas-grants", "ASAPRegister"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/ideas-grants#Funding%20Partners", "ASAPRegister"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-gacd-funding-call", "ASAPRegister"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc", "Accepting"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/manage-your-grant", "Accepting"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-amed-2024-aspire", "Adopting"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/2024-nhmrc-cihr-ccna", "Aging"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/postgraduate-scholarships", "April"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/nhmrc-collaborations-in-health-services-research", "August"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls", "Australian"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/cre", "aims"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/synergy-grants", "and"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/clinical-trials-and-cohort-studies-grants", "and"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/investigator-grants", "approvals,"),
    @("https://sites.research.unimelb.edu.au/research-funding/nhmrc/targeted-calls/eucrg", "available4")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
